# Weekly update: insert a new Choclo price-report row for
# "Vega Monumental Concepción" as row 106, pushing the existing
# rows 106-127 down to 107-128.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 106 (shifts rows 106..127 -> 107..128)
$ws.Rows.Item(106).Insert()

# Populate the newly inserted row with this week's data
$ws.Cells.Item(106, 1).Value = 11
$ws.Cells.Item(106, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(106, 3).Value = "Bíobío"
$ws.Cells.Item(106, 4).Value = 44943
$ws.Cells.Item(106, 5).Value = 8
$ws.Cells.Item(106, 6).Value = 100112024
$ws.Cells.Item(106, 7).Value = "Choclo"
$ws.Cells.Item(106, 8).Value = "Choclero"
$ws.Cells.Item(106, 9).Value = "Primera"
$ws.Cells.Item(106, 10).Value = 7000
$ws.Cells.Item(106, 11).Value = 200
$ws.Cells.Item(106, 12).Value = 250
$ws.Cells.Item(106, 13).Value = 236
$ws.Cells.Item(106, 14).Value = "$/unidad"
$ws.Cells.Item(106, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(106, 16).Value = 236
$ws.Cells.Item(106, 17).Value = 1
$ws.Cells.Item(106, 18).Value = "Hortaliza"
